$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before
#    the current second sheet, which is "2022-Q2").
# ---------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$ws = $wb.Worksheets.Add($refSheet)
$ws.Name = "2022-Q3"

# Header row.
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data row for the single fund holding in 2022-Q3.
$ws.Range("A2").Value = 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "512290"
$ws.Range("C2").Value = "国泰中证生物医药ETF"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "99.74"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2.15"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.8750"
$ws.Range("H2").Value = 8

# Match the bold / thin-border / centered-top style that every other
# quarter sheet uses for its header row and row-index column, by
# copying it (format only) from the neighbouring "2022-Q2" sheet.
$styleSrc = $wb.Worksheets.Item("2022-Q2")
$styleSrc.Range("B1:H1").Copy() | Out-Null
$ws.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$styleSrc.Range("A2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    above the existing data (shifting the rest down, which is
#    handled automatically since we overwrite every row's contents).
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$dates = @("2022-Q3", "2022-Q2", "2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$counts = @(1, 2, 6, 13, 3, 24, 6, 2)
$values = @(0.88, 1.14, 1.08, 17.18, 2.49, 37.71, 3.12, 0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $dates[$i]
    $summary.Cells.Item($row, 3).Value = $counts[$i]
    $summary.Cells.Item($row, 4).Value = $values[$i]
}
